$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.158.90'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -4.86%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.233.56'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -5.58%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.22%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '318.12'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.51%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '99.71'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -9.44%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.580'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -7.73%  '

# Row 8
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.15%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.566'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -8.36%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.65'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -11.35%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '54.30'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -2.10%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0819'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -10.85%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.70'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -9.28%  '

# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.29%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.574.28'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -5.43%  '

# Row 16
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.865'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -12.08%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.14'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -8.08%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.229.55'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -5.49%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '43.160.00'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -4.74%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.39'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +4.90%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0969'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -9.16%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.52'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -10.93%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.40'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -10.76%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.17'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -7.90%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '235.65'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -9.03%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.17'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -5.62%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.17%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.20'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -8.55%  '

# Row 29
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -7.31%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.41'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -13.27%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0889'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -8.64%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.59'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -8.19%  '

# Row 33
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '158.02'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -7.32%  '

# Row 34
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '33.87'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -11.26%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.76'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -5.63%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.33'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +12.16%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.02'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +16.01%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.122'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -6.78%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.49'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -7.09%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.105'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -8.10%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -8.56%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0324'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -9.26%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.01'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.72%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.815.69'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +8.87%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '12.05'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -6.13%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '88.11'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -10.97%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.81'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +67.58%  '

# Row 48
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.55'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.69%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.207'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -11.05%  '

# Row 50
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '77.41'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -6.25%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '60.78'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -13.12%  '
